$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Myoc"
$ws.Cells.Item(2,3).Value = "Fzd4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.714093
$ws.Cells.Item(2,8).Value = 2.142279
$ws.Cells.Item(2,9).Value = 0.001553869939232348
$ws.Cells.Item(2,10).Value = 0.001553869939232348
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 13.57958433333333
$ws.Cells.Item(2,14).Value = 40.738753
$ws.Cells.Item(2,15).Value = 0.2289698008477291
$ws.Cells.Item(2,16).Value = 0.2289698008477291
$ws.Cells.Item(2,17).Value = 9.697086115343001
$ws.Cells.Item(2,18).Value = 87.273775038087
$ws.Cells.Item(2,19).Value = 0.0003557892905293035
$ws.Cells.Item(2,20).Value = 0.0003557892905293035
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Myoc"
$ws.Cells.Item(3,3).Value = "Fzd4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.714093
$ws.Cells.Item(3,8).Value = 2.142279
$ws.Cells.Item(3,9).Value = 0.001553869939232348
$ws.Cells.Item(3,10).Value = 0.001553869939232348
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 19.768727
$ws.Cells.Item(3,14).Value = 59.306181
$ws.Cells.Item(3,15).Value = 0.3333269541315948
$ws.Cells.Item(3,16).Value = 0.3333269541315948
$ws.Cells.Item(3,17).Value = 14.116709569611
$ws.Cells.Item(3,18).Value = 127.050386126499
$ws.Cells.Item(3,19).Value = 0.0005179467339609647
$ws.Cells.Item(3,20).Value = 0.0005179467339609647
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Myoc"
$ws.Cells.Item(4,3).Value = "Fzd4"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.714093
$ws.Cells.Item(4,8).Value = 2.142279
$ws.Cells.Item(4,9).Value = 0.001553869939232348
$ws.Cells.Item(4,10).Value = 0.001553869939232348
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 25.95900466666667
$ws.Cells.Item(4,14).Value = 77.877014
$ws.Cells.Item(4,15).Value = 0.4377032450206762
$ws.Cells.Item(4,16).Value = 0.4377032450206762
$ws.Cells.Item(4,17).Value = 18.537143519434
$ws.Cells.Item(4,18).Value = 166.834291674906
$ws.Cells.Item(4,19).Value = 0.0006801339147420796
$ws.Cells.Item(4,20).Value = 0.0006801339147420796
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Myoc"
$ws.Cells.Item(5,3).Value = "Fzd4"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 453.4108886666666
$ws.Cells.Item(5,8).Value = 1360.232666
$ws.Cells.Item(5,9).Value = 0.9866243612803347
$ws.Cells.Item(5,10).Value = 0.9866243612803348
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 13.57958433333333
$ws.Cells.Item(5,14).Value = 40.738753
$ws.Cells.Item(5,15).Value = 0.2289698008477291
$ws.Cells.Item(5,16).Value = 0.2289698008477291
$ws.Cells.Item(5,17).Value = 6157.131400300611
$ws.Cells.Item(5,18).Value = 55414.1826027055
$ws.Cells.Item(5,19).Value = 0.2259071835138761
$ws.Cells.Item(5,20).Value = 0.2259071835138761
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Myoc"
$ws.Cells.Item(6,3).Value = "Fzd4"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 453.4108886666666
$ws.Cells.Item(6,8).Value = 1360.232666
$ws.Cells.Item(6,9).Value = 0.9866243612803347
$ws.Cells.Item(6,10).Value = 0.9866243612803348
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 19.768727
$ws.Cells.Item(6,14).Value = 59.306181
$ws.Cells.Item(6,15).Value = 0.3333269541315948
$ws.Cells.Item(6,16).Value = 0.3333269541315948
$ws.Cells.Item(6,17).Value = 8963.356076878725
$ws.Cells.Item(6,18).Value = 80670.20469190853
$ws.Cells.Item(6,19).Value = 0.3288684932176041
$ws.Cells.Item(6,20).Value = 0.3288684932176042
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Myoc"
$ws.Cells.Item(7,3).Value = "Fzd4"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 453.4108886666666
$ws.Cells.Item(7,8).Value = 1360.232666
$ws.Cells.Item(7,9).Value = 0.9866243612803347
$ws.Cells.Item(7,10).Value = 0.9866243612803348
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 25.95900466666667
$ws.Cells.Item(7,14).Value = 77.877014
$ws.Cells.Item(7,15).Value = 0.4377032450206762
$ws.Cells.Item(7,16).Value = 0.4377032450206762
$ws.Cells.Item(7,17).Value = 11770.09537481548
$ws.Cells.Item(7,18).Value = 105930.8583733393
$ws.Cells.Item(7,19).Value = 0.4318486845488545
$ws.Cells.Item(7,20).Value = 0.4318486845488546
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Myoc"
$ws.Cells.Item(8,3).Value = "Fzd4"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 5.432785666666667
$ws.Cells.Item(8,8).Value = 16.298357
$ws.Cells.Item(8,9).Value = 0.01182176878043295
$ws.Cells.Item(8,10).Value = 0.01182176878043295
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 13.57958433333333
$ws.Cells.Item(8,14).Value = 40.738753
$ws.Cells.Item(8,15).Value = 0.2289698008477291
$ws.Cells.Item(8,16).Value = 0.2289698008477291
$ws.Cells.Item(8,17).Value = 73.77497112542457
$ws.Cells.Item(8,18).Value = 663.974740128821
$ws.Cells.Item(8,19).Value = 0.002706828043323633
$ws.Cells.Item(8,20).Value = 0.002706828043323633
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Myoc"
$ws.Cells.Item(9,3).Value = "Fzd4"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 5.432785666666667
$ws.Cells.Item(9,8).Value = 16.298357
$ws.Cells.Item(9,9).Value = 0.01182176878043295
$ws.Cells.Item(9,10).Value = 0.01182176878043295
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 19.768727
$ws.Cells.Item(9,14).Value = 59.306181
$ws.Cells.Item(9,15).Value = 0.3333269541315948
$ws.Cells.Item(9,16).Value = 0.3333269541315948
$ws.Cells.Item(9,17).Value = 107.3992566938463
$ws.Cells.Item(9,18).Value = 966.5933102446169
$ws.Cells.Item(9,19).Value = 0.003940514180029692
$ws.Cells.Item(9,20).Value = 0.003940514180029692
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Myoc"
$ws.Cells.Item(10,3).Value = "Fzd4"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 5.432785666666667
$ws.Cells.Item(10,8).Value = 16.298357
$ws.Cells.Item(10,9).Value = 0.01182176878043295
$ws.Cells.Item(10,10).Value = 0.01182176878043295
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 25.95900466666667
$ws.Cells.Item(10,14).Value = 77.877014
$ws.Cells.Item(10,15).Value = 0.4377032450206762
$ws.Cells.Item(10,16).Value = 0.4377032450206762
$ws.Cells.Item(10,17).Value = 141.0297084739998
$ws.Cells.Item(10,18).Value = 1269.267376265998
$ws.Cells.Item(10,19).Value = 0.005174426557079622
$ws.Cells.Item(10,20).Value = 0.005174426557079622
